$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L),
# bringing in a new fiscal-year column of data.
$ws.Columns("D").Insert()

# Copy the formatting (number format/style) from column E (the old D,
# now shifted right) into the freshly inserted blank column D so the
# new cells match the style of the data they precede. Done per data
# block so we don't create stray cells on the section-header rows
# (5, 6, 37, 79) that never had data in columns D:K.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal year figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1285600
$ws.Range("D9").Value = 908500
$ws.Range("D10").Value = 377200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 936100
$ws.Range("D18").Value = 349600
$ws.Range("D20").Value = 2500
$ws.Range("D21").Value = 460900
$ws.Range("D22").Value = 21400
$ws.Range("D23").Value = 330600
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 330600
$ws.Range("D27").Value = 49300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2500
$ws.Range("D33").Value = 49300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 49300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 19700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 222600
$ws.Range("D44").Value = 4400
$ws.Range("D45").Value = 400
$ws.Range("D46").Value = 247100
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 1148300
$ws.Range("D49").Value = 5800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 60900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1462100
$ws.Range("D57").Value = 32500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 16300
$ws.Range("D60").Value = 48800
$ws.Range("D61").Value = 477600
$ws.Range("D62").Value = 1700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1246300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = -242600
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 458400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 49300
$ws.Range("D83").Value = 108800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 436200
$ws.Range("D91").Value = -39900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -51800
$ws.Range("D96").Value = -53400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -391600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -7300
